$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.506.56'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.622.41'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E4').Value = '  +0.00%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '211.68'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.56%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.520'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('E7').Value = '  -0.02%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '23.18'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.21%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.263'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('E10').Value = '  +0.01%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0882'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('D12').Value = '1.852.53'
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('D13').Value = '1.633.07'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('E14').Value = '  -0.19%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.551'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.99%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '65.36'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').Value = '27.488.88'
$ws.Range('E17').Value = '  -0.58%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '229.43'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').Value = '0.0₃0719'
$ws.Range('E19').Value = '  -0.66%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.55'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('E21').Value = '  +0.01%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '10.44'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +4.03%  '
$ws.Range('E23').Value = '  +1.23%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '149.27'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('E28').Value = '  -0.02%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '15.51'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  -0.69%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.0484'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('D33').Value = '1.466.62'
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('E34').Value = '  -2.36%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.55'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -1.61%  '
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('E37').Value = '  +4.85%  '
$ws.Range('E38').Value = '  +0.17%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.874'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.58%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.553'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -2.52%  '
$ws.Range('E41').Value = '  -0.44%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '68.03'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -3.21%  '
$ws.Range('E44').Value = '  +0.85%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.20'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -1.99%  '
$ws.Range('E46').Value = '  -4.88%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.75'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.23%  '
$ws.Range('D48').Value = '1.762.31'
$ws.Range('E48').Value = '  -1.13%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '87.19'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.18%  '
$ws.Range('E50').Value = '  -0.93%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0994'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.49%  '
